# Auto-generated edit script: update Leve profit-calc columns (H-N) per refreshed market prices
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 6420.9414
$ws.Range("I33").Value = 57.125
$ws.Range("J33").Value = 12077.667
$ws.Range("K33").Value = 57.125
$ws.Range("L33").Value = 12077.667
$ws.Range("M33").Value = 171.875
$ws.Range("N33").Value = -12535.667
$ws.Range("H113").Value = 2037.475
$ws.Range("I113").Value = 1909.2258
$ws.Range("J113").Value = 2479.2222
$ws.Range("K113").Value = 1909.2258
$ws.Range("L113").Value = 2479.2222
$ws.Range("M113").Value = 1344.7742
$ws.Range("N113").Value = -8987.2222
$ws.Range("H130").Value = 36291.11
$ws.Range("J130").Value = 36291.11
$ws.Range("L130").Value = 36291.11
$ws.Range("N130").Value = -46331.11
$ws.Range("H135").Value = 1650.0731
$ws.Range("I135").Value = 1034.7715
$ws.Range("J135").Value = 5239.3335
$ws.Range("K135").Value = 9312.943500000001
$ws.Range("L135").Value = 47154.0015
$ws.Range("M135").Value = -6777.943500000001
$ws.Range("N135").Value = -52224.0015
$ws.Range("H137").Value = 985.09375
$ws.Range("I137").Value = 853.62964
$ws.Range("J137").Value = 1695
$ws.Range("K137").Value = 2560.88892
$ws.Range("L137").Value = 5085
$ws.Range("M137").Value = -10.88891999999987
$ws.Range("N137").Value = -10185

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1600.84
$ws.Range("I2").Value = 729.1177
$ws.Range("J2").Value = 3453.25
$ws.Range("K2").Value = 729.1177
$ws.Range("L2").Value = 3453.25
$ws.Range("M2").Value = -616.1177
$ws.Range("N2").Value = -3679.25
$ws.Range("H32").Value = 342935.94
$ws.Range("I32").Value = 2369.3418
$ws.Range("J32").Value = 5723888.5
$ws.Range("K32").Value = 2369.3418
$ws.Range("L32").Value = 5723888.5
$ws.Range("M32").Value = -2082.3418
$ws.Range("N32").Value = -5724462.5
$ws.Range("H74").Value = 1224.8334
$ws.Range("I74").Value = 1224.75
$ws.Range("J74").Value = 1225
$ws.Range("K74").Value = 1224.75
$ws.Range("L74").Value = 1225
$ws.Range("M74").Value = -350.75
$ws.Range("N74").Value = -2973
$ws.Range("H77").Value = 1224.8334
$ws.Range("I77").Value = 1224.75
$ws.Range("J77").Value = 1225
$ws.Range("K77").Value = 6123.75
$ws.Range("L77").Value = 6125
$ws.Range("M77").Value = -1755.75
$ws.Range("N77").Value = -14861
$ws.Range("H116").Value = 1600.84
$ws.Range("I116").Value = 729.1177
$ws.Range("J116").Value = 3453.25
$ws.Range("K116").Value = 729.1177
$ws.Range("L116").Value = 3453.25
$ws.Range("M116").Value = 1564.8823
$ws.Range("N116").Value = -8041.25
$ws.Range("H122").Value = 14314.388
$ws.Range("I122").Value = 17136.574
$ws.Range("J122").Value = 1771.3334
$ws.Range("K122").Value = 51409.722
$ws.Range("L122").Value = 5314.0002
$ws.Range("M122").Value = -48959.722
$ws.Range("N122").Value = -10214.0002
$ws.Range("H132").Value = 22730998
$ws.Range("I132").Value = 31250990
$ws.Range("J132").Value = 11019.5
$ws.Range("K132").Value = 93752970
$ws.Range("L132").Value = 33058.5
$ws.Range("M132").Value = -93750440
$ws.Range("N132").Value = -38118.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1600.84
$ws.Range("I3").Value = 729.1177
$ws.Range("J3").Value = 3453.25
$ws.Range("K3").Value = 729.1177
$ws.Range("L3").Value = 3453.25
$ws.Range("M3").Value = -615.1177
$ws.Range("N3").Value = -3681.25
$ws.Range("H99").Value = 1423.3334
$ws.Range("I99").Value = 803.5454999999999
$ws.Range("J99").Value = 3127.75
$ws.Range("K99").Value = 803.5454999999999
$ws.Range("L99").Value = 3127.75
$ws.Range("M99").Value = 694.4545000000001
$ws.Range("N99").Value = -6123.75
$ws.Range("H134").Value = 7730.074
$ws.Range("I134").Value = 2596.48
$ws.Range("J134").Value = 71900
$ws.Range("K134").Value = 7789.440000000001
$ws.Range("L134").Value = 215700
$ws.Range("M134").Value = -5254.440000000001
$ws.Range("N134").Value = -220770

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2222.625
$ws.Range("I31").Value = 1530.1666
$ws.Range("J31").Value = 4300
$ws.Range("K31").Value = 1530.1666
$ws.Range("L31").Value = 4300
$ws.Range("M31").Value = -1235.1666
$ws.Range("N31").Value = -4890
$ws.Range("H34").Value = 2222.625
$ws.Range("I34").Value = 1530.1666
$ws.Range("J34").Value = 4300
$ws.Range("K34").Value = 1530.1666
$ws.Range("L34").Value = 4300
$ws.Range("M34").Value = -1328.1666
$ws.Range("N34").Value = -4704
$ws.Range("H58").Value = 1146.6957
$ws.Range("I58").Value = 857.375
$ws.Range("J58").Value = 1808
$ws.Range("K58").Value = 857.375
$ws.Range("L58").Value = 1808
$ws.Range("M58").Value = -654.375
$ws.Range("N58").Value = -2214
$ws.Range("H132").Value = 34798.566
$ws.Range("I132").Value = 621.8095
$ws.Range("J132").Value = 114544.336
$ws.Range("K132").Value = 1865.4285
$ws.Range("L132").Value = 343633.008
$ws.Range("M132").Value = 664.5715
$ws.Range("N132").Value = -348693.008
$ws.Range("H136").Value = 1146.6957
$ws.Range("I136").Value = 857.375
$ws.Range("J136").Value = 1808
$ws.Range("K136").Value = 2572.125
$ws.Range("L136").Value = 5424
$ws.Range("M136").Value = -22.125
$ws.Range("N136").Value = -10524

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 30000.4
$ws.Range("J37").Value = 30000.4
$ws.Range("L37").Value = 90001.20000000001
$ws.Range("N37").Value = -90225.20000000001
$ws.Range("H75").Value = 1654.6666
$ws.Range("I75").Value = 804.3333
$ws.Range("J75").Value = 2505
$ws.Range("K75").Value = 2412.9999
$ws.Range("L75").Value = 7515
$ws.Range("M75").Value = -1414.9999
$ws.Range("N75").Value = -9511
$ws.Range("H78").Value = 1654.6666
$ws.Range("I78").Value = 804.3333
$ws.Range("J78").Value = 2505
$ws.Range("K78").Value = 7238.9997
$ws.Range("L78").Value = 22545
$ws.Range("M78").Value = -2246.9997
$ws.Range("N78").Value = -32529
$ws.Range("H131").Value = 6098542
$ws.Range("J131").Value = 7247396
$ws.Range("L131").Value = 21742188
$ws.Range("N131").Value = -21752268
$ws.Range("H132").Value = 1807
$ws.Range("I132").Value = 1113.8572
$ws.Range("J132").Value = 2500.1428
$ws.Range("K132").Value = 10024.7148
$ws.Range("L132").Value = 22501.2852
$ws.Range("M132").Value = -7494.7148
$ws.Range("N132").Value = -27561.2852
$ws.Range("H137").Value = 2557.9033
$ws.Range("I137").Value = 932.8570999999999
$ws.Range("J137").Value = 3031.875
$ws.Range("K137").Value = 2798.5713
$ws.Range("L137").Value = 9095.625
$ws.Range("M137").Value = 2301.4287
$ws.Range("N137").Value = -19295.625

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2033.3422
$ws.Range("I68").Value = 1898.75
$ws.Range("K68").Value = 1898.75
$ws.Range("M68").Value = -1149.75
$ws.Range("H71").Value = 2033.3422
$ws.Range("I71").Value = 1898.75
$ws.Range("K71").Value = 9493.75
$ws.Range("M71").Value = -5749.75
$ws.Range("H132").Value = 3937.1206
$ws.Range("I132").Value = 5116.8945
$ws.Range("J132").Value = 1695.55
$ws.Range("K132").Value = 15350.6835
$ws.Range("L132").Value = 5086.65
$ws.Range("M132").Value = -12820.6835
$ws.Range("N132").Value = -10146.65

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 16875514
$ws.Range("I132").Value = 23150828
$ws.Range("J132").Value = 1472468.2
$ws.Range("K132").Value = 69452484
$ws.Range("L132").Value = 4417404.6
$ws.Range("M132").Value = -69449954
$ws.Range("N132").Value = -4422464.6
